$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C9").Value = "01/06/2016 (Olivier)"
